$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly crypto price/volume refresh (GitHub Actions cron).
# D-column numeric-looking quotes must stay text -> force quote-prefix via
# NumberFormat "@", then restore the default "Normal" style so no new
# cell-level style index is introduced (matches source which had no 's' attr).

$ws.Range("D2").Value = '26.856.19'
$ws.Range("E2").Value = '  -0.92%  '

$ws.Range("D3").Value = '1.856.34'
$ws.Range("E3").Value = '  -0.59%  '

$ws.Range("E4").Value = '  -0.19%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '304.87'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.75%  '

$ws.Range("E6").Value = '  -0.18%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5071'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.28%  '

$ws.Range("E8").Value = '  -2.69%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07177'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.25%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8897'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.67%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.63'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.36%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07528'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.44%  '

$ws.Range("D13").Value = '1.845.92'
$ws.Range("E13").Value = '  -1.27%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '91.70'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.07%  '

$ws.Range("E15").Value = '  -1.80%  '

$ws.Range("E16").Value = '  -0.18%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008531'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.16%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '14.06'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.66%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.0000'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.10%  '

$ws.Range("D20").Value = '26.908.15'
$ws.Range("E20").Value = '  -0.90%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.023'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.08%  '

$ws.Range("D22").Value = '2.084.89'
$ws.Range("E22").Value = '  -1.63%  '

$ws.Range("E23").Value = '  -2.89%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.450'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.26%  '

$ws.Range("B25").Value = 'Toncoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.808'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.78%  '

$ws.Range("B26").Value = 'Monero'
$ws.Range("C26").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '146.10'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.50%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.82'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.06%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.051'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -5.06%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '112.95'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.07%  '

$ws.Range("E30").Value = '  -2.08%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.663'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.51%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09210'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.21%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05096'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.91%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.975'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -3.88%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7340'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.73%  '

$ws.Range("E36").Value = '  -1.96%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.216'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +5.96%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02003'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.75%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.481'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.32%  '

$ws.Range("E40").Value = '  -0.90%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.5321'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.40%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '118.56'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +3.69%  '

$ws.Range("E43").Value = '  -2.29%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.351'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.88%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.1471'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.89%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4633'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.05%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.9999'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.17%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.950'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.41%  '

$ws.Range("E49").Value = '  -0.82%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '36.97'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.25%  '

$ws.Range("E51").Value = '  -3.18%  '
